$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.255028963088989
$ws.Range("B1").Value = 2.005975961685181
$ws.Range("C1").Value = 5.812239646911621
$ws.Range("D1").Value = 1.952891111373901
$ws.Range("E1").Value = 1.131213665008545
